# Apply updated cryptocurrency price/volume figures (columns D and E).
# Source cells store plain display text (e.g. "25.720.13", "  -1.22%  "),
# not numbers, so each write forces Text formatting, assigns the literal
# string, then restores General/Normal formatting to avoid leaving any
# lingering cell-format differences behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "25.720.13" },
    @{ Cell = "E2"; Value = "  -1.22%  " },
    @{ Cell = "D3"; Value = "1.622.19" },
    @{ Cell = "E3"; Value = "  -1.17%  " },
    @{ Cell = "E4"; Value = "  -0.12%  " },
    @{ Cell = "D5"; Value = "214.49" },
    @{ Cell = "E5"; Value = "  -0.57%  " },
    @{ Cell = "D6"; Value = "0.5075" },
    @{ Cell = "E6"; Value = "  -1.12%  " },
    @{ Cell = "D7"; Value = "1.002" },
    @{ Cell = "E7"; Value = "  -0.08%  " },
    @{ Cell = "D8"; Value = "0.2560" },
    @{ Cell = "E8"; Value = "  -1.22%  " },
    @{ Cell = "D9"; Value = "0.06366" },
    @{ Cell = "E9"; Value = "  -0.26%  " },
    @{ Cell = "E10"; Value = "  -3.29%  " },
    @{ Cell = "D11"; Value = "0.07761" },
    @{ Cell = "E11"; Value = "  -0.22%  " },
    @{ Cell = "D12"; Value = "4.234" },
    @{ Cell = "E12"; Value = "  -1.71%  " },
    @{ Cell = "D13"; Value = "1.623.58" },
    @{ Cell = "E13"; Value = "  -0.68%  " },
    @{ Cell = "D14"; Value = "1.845.11" },
    @{ Cell = "E14"; Value = "  -1.20%  " },
    @{ Cell = "D15"; Value = "0.5542" },
    @{ Cell = "E15"; Value = "  +1.03%  " },
    @{ Cell = "D16"; Value = "63.47" },
    @{ Cell = "E16"; Value = "  -1.90%  " },
    @{ Cell = "D17"; Value = "0.0₅7528" },
    @{ Cell = "E17"; Value = "  -3.09%  " },
    @{ Cell = "D18"; Value = "25.740.29" },
    @{ Cell = "D19"; Value = "1.002" },
    @{ Cell = "E19"; Value = "  -0.07%  " },
    @{ Cell = "D20"; Value = "193.45" },
    @{ Cell = "E20"; Value = "  -2.90%  " },
    @{ Cell = "D21"; Value = "4.376" },
    @{ Cell = "E21"; Value = "  -1.83%  " },
    @{ Cell = "D22"; Value = "9.745" },
    @{ Cell = "E22"; Value = "  -2.59%  " },
    @{ Cell = "D23"; Value = "5.946" },
    @{ Cell = "E23"; Value = "  -2.58%  " },
    @{ Cell = "E24"; Value = "  -0.16%  " },
    @{ Cell = "D25"; Value = "1.857" },
    @{ Cell = "E25"; Value = "  -1.87%  " },
    @{ Cell = "D26"; Value = "140.50" },
    @{ Cell = "E26"; Value = "  -1.32%  " },
    @{ Cell = "E27"; Value = "  +1.28%  " },
    @{ Cell = "D28"; Value = "6.718" },
    @{ Cell = "E28"; Value = "  -2.50%  " },
    @{ Cell = "D29"; Value = "15.46" },
    @{ Cell = "E29"; Value = "  -1.40%  " },
    @{ Cell = "D30"; Value = "1.233" },
    @{ Cell = "E30"; Value = "  -0.51%  " },
    @{ Cell = "D31"; Value = "0.04860" },
    @{ Cell = "E31"; Value = "  -1.50%  " },
    @{ Cell = "D32"; Value = "3.292" },
    @{ Cell = "E32"; Value = "  -0.51%  " },
    @{ Cell = "D33"; Value = "3.168" },
    @{ Cell = "E33"; Value = "  -1.69%  " },
    @{ Cell = "D34"; Value = "1.540" },
    @{ Cell = "E34"; Value = "  -0.51%  " },
    @{ Cell = "E35"; Value = "  -0.60%  " },
    @{ Cell = "E36"; Value = "  -3.60%  " },
    @{ Cell = "D37"; Value = "1.124.36" },
    @{ Cell = "E37"; Value = "  +0.82%  " },
    @{ Cell = "D38"; Value = "2.533" },
    @{ Cell = "E38"; Value = "  -2.12%  " },
    @{ Cell = "D39"; Value = "0.5474" },
    @{ Cell = "E39"; Value = "  -2.05%  " },
    @{ Cell = "D40"; Value = "0.01556" },
    @{ Cell = "E40"; Value = "  -1.09%  " },
    @{ Cell = "D41"; Value = "1.001" },
    @{ Cell = "E41"; Value = "  -0.13%  " },
    @{ Cell = "D42"; Value = "5.567" },
    @{ Cell = "E42"; Value = "  +0.35%  " },
    @{ Cell = "D43"; Value = "0.7929" },
    @{ Cell = "E43"; Value = "  -2.39%  " },
    @{ Cell = "D44"; Value = "96.99" },
    @{ Cell = "E44"; Value = "  -2.89%  " },
    @{ Cell = "D45"; Value = "1.769.81" },
    @{ Cell = "E45"; Value = "  -0.48%  " },
    @{ Cell = "E46"; Value = "  -7.92%  " },
    @{ Cell = "D47"; Value = "0.4414" },
    @{ Cell = "E47"; Value = "  -2.62%  " },
    @{ Cell = "D48"; Value = "54.54" },
    @{ Cell = "E48"; Value = "  -1.39%  " },
    @{ Cell = "D49"; Value = "0.05119" },
    @{ Cell = "E49"; Value = "  -3.08%  " },
    @{ Cell = "D50"; Value = "7.557" },
    @{ Cell = "E50"; Value = "  +2.44%  " },
    @{ Cell = "D51"; Value = "0.9963" },
    @{ Cell = "E51"; Value = "  -1.09%  " }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    $r.NumberFormat = "@"
    $r.Value = $u.Value
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}
